$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column before D, shifting existing D:L data right by one column.
$ws.Columns("D:D").Insert(-4161)

# Copy number formats / styles from the (now shifted) column E into the
# freshly inserted column D so every row keeps its original per-row style
# (date style for header rows, numeric style for data rows).
$ws.Range("E7:L102").Copy()
$ws.Range("D7:D102").PasteSpecial(-4122)

# Populate the new column D with the latest reported quarter's figures.
# Period-ending header rows (new quarter date, serial 43373 = 2018-09-07).
$ws.Range("D7").Value = 43373
$ws.Range("D38").Value = 43373
$ws.Range("D80").Value = 43373

# Income statement
$ws.Range("D8").Value = 5500
$ws.Range("D9").Value = "NA"
$ws.Range("D10").Value = "NA"
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("D17").Value = 1000
$ws.Range("D18").Value = 4500
$ws.Range("D20").Value = -3000
$ws.Range("D21").Value = 1600
$ws.Range("D22").Value = 0
$ws.Range("D23").Value = 1400
$ws.Range("D24").Value = 300
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 1100
$ws.Range("D27").Value = 1100
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 3000
$ws.Range("D33").Value = 1100
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 1100

# Balance sheet
$ws.Range("D41").Value = 14900
$ws.Range("D42").Value = 500
$ws.Range("D43").Value = 0
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 0
$ws.Range("D46").Value = 0
$ws.Range("D47").Value = 0
$ws.Range("D48").Value = 12600
$ws.Range("D49").Value = 0
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 600
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 512300
$ws.Range("D57").Value = 200
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 0
$ws.Range("D60").Value = 0
$ws.Range("D61").Value = 14300
$ws.Range("D62").Value = 0
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 476500
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = -22600
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 35800
$ws.Range("D77").Value = 0

# Cash flow statement
$ws.Range("D81").Value = 1100
$ws.Range("D83").Value = 200
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 5600
$ws.Range("D91").Value = 0
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = -11200
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = 5200
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = -300
